# Updates the cryptos list (Price column D, Volume(1h) column E)
# All target cells are stored as text in the workbook, so Price cells
# are explicitly written with a text number format to stop Excel from
# auto-converting numeric-looking strings (e.g. "148.60") into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '64.802.35'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.146.04'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '575.74'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '148.60'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.73%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.147.02'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -0.55%  '
$ws.Range('E10').Value = '  -3.02%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.07'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.92%  '
$ws.Range('E12').Value = '  -1.63%  '
$ws.Range('E13').Value = '  +0.70%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '36.93'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.662.20'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.932.35'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.148.91'
$ws.Range('D17').Style = "Normal"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.07'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.85%  '
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '502.22'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.56%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.74'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.15%  '
$ws.Range('E22').Value = '  -3.28%  '
$ws.Range('E23').Value = '  -2.29%  '
$ws.Range('E24').Value = '  -2.32%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '83.67'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.30%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.997'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('E27').Value = '  -1.65%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.81'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.16%  '
$ws.Range('E29').Value = '  -1.32%  '
$ws.Range('E30').Value = '  +5.53%  '
$ws.Range('E31').Value = '  -2.06%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.999'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.16'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.08%  '
$ws.Range('E35').Value = '  -2.57%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '54.51'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.05%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0887'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +3.58%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '475.23'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('E39').Value = '  -2.58%  '
$ws.Range('E40').Value = '  -3.94%  '
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.002.90'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.67%  '
$ws.Range('E43').Value = '  -3.43%  '
$ws.Range('E44').Value = '  -3.79%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.40'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.05%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '27.98'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.37%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0₃0577'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.38%  '
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('E49').Value = '  -2.13%  '
$ws.Range('E50').Value = '  -3.55%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '117.15'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.15%  '
